# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The underlying macro-generated report re-sorted / re-generated the detail
# rows (16-28) of the "Estado de Cuenta" table: a new worker record
# (73157992 - SALVADOR FRIERI DEL CASTILLO, period 1607) was moved to the
# top of the table (row 16) with updated Valor Mora / Salario Basico values,
# and the remaining periods for MARY TORRES RIPOLL (33333269) and
# BLANCA ROSA TORRES MUNOZ (23002667) were re-interleaved with refreshed
# Salario Basico figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: SALVADOR FRIERI DEL CASTILLO - 1607
$ws.Range("C16").Value = "73157992"
$ws.Range("D16").Value = "SALVADOR FRIERI DEL CASTILLO"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: MARY TORRES RIPOLL - 1801
$ws.Range("C17").Value = "33333269"
$ws.Range("D17").Value = "MARY TORRES RIPOLL"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242

# Row 18: MARY TORRES RIPOLL - 1803
$ws.Range("C18").Value = "33333269"
$ws.Range("D18").Value = "MARY TORRES RIPOLL"
$ws.Range("E18").Value = "1803"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242

# Row 19: MARY TORRES RIPOLL - 1804
$ws.Range("C19").Value = "33333269"
$ws.Range("D19").Value = "MARY TORRES RIPOLL"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242

# Row 20: MARY TORRES RIPOLL - 1805
$ws.Range("C20").Value = "33333269"
$ws.Range("D20").Value = "MARY TORRES RIPOLL"
$ws.Range("E20").Value = "1805"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242

# Row 21: MARY TORRES RIPOLL - 1806
$ws.Range("C21").Value = "33333269"
$ws.Range("D21").Value = "MARY TORRES RIPOLL"
$ws.Range("E21").Value = "1806"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 781242

# Row 22: BLANCA ROSA TORRES MUNOZ - 1806
$ws.Range("C22").Value = "23002667"
$ws.Range("D22").Value = "BLANCA ROSA TORRES MUNOZ"
$ws.Range("E22").Value = "1806"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 781242

# Row 23: MARY TORRES RIPOLL - 1807
$ws.Range("C23").Value = "33333269"
$ws.Range("D23").Value = "MARY TORRES RIPOLL"
$ws.Range("E23").Value = "1807"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 781242

# Row 24: BLANCA ROSA TORRES MUNOZ - 1807
$ws.Range("C24").Value = "23002667"
$ws.Range("D24").Value = "BLANCA ROSA TORRES MUNOZ"
$ws.Range("E24").Value = "1807"
$ws.Range("F24").Value = 31249
$ws.Range("G24").Value = 781242

# Row 25: MARY TORRES RIPOLL - 1808
$ws.Range("C25").Value = "33333269"
$ws.Range("D25").Value = "MARY TORRES RIPOLL"
$ws.Range("E25").Value = "1808"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 781242

# Row 26: BLANCA ROSA TORRES MUNOZ - 1808
$ws.Range("C26").Value = "23002667"
$ws.Range("D26").Value = "BLANCA ROSA TORRES MUNOZ"
$ws.Range("E26").Value = "1808"
$ws.Range("F26").Value = 31249
$ws.Range("G26").Value = 781242

# Row 27: MARY TORRES RIPOLL - 1809
$ws.Range("C27").Value = "33333269"
$ws.Range("D27").Value = "MARY TORRES RIPOLL"
$ws.Range("E27").Value = "1809"
$ws.Range("F27").Value = 31249
$ws.Range("G27").Value = 781242

# Row 28: BLANCA ROSA TORRES MUNOZ - 1809
$ws.Range("C28").Value = "23002667"
$ws.Range("D28").Value = "BLANCA ROSA TORRES MUNOZ"
$ws.Range("E28").Value = "1809"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242
